# Insert two new weekly price records (Primera / Segunda quality, date 44551)
# right before the existing row 344, pushing the rest of the table down by
# two rows. This matches the diff: dimension grows from A1:T446 to A1:T448
# and all former rows 344-446 become rows 346-448 unchanged.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 2 blank rows above the current row 344.
$ws.Rows("344:345").Insert()

# --- New row 344: Primera, $/caja 12 unidades ---
$ws.Range("A344").Value = 3
$ws.Range("B344").Value = "Femacal de La Calera"
$ws.Range("C344").Value = "Coquimbo"
$ws.Range("D344").Value = 44551
$ws.Range("D344").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("E344").Value = 5
$ws.Range("F344").Value = "Fruta"
$ws.Range("G344").Value = 100108
$ws.Range("H344").Value = "Tropicales y subtropicales"
$ws.Range("I344").Value = 100108005
$ws.Range("J344").Value = "Piña"
$ws.Range("K344").Value = "Caramelo"
$ws.Range("L344").Value = "Primera"
$ws.Range("M344").Value = 162
$ws.Range("N344").Value = 18000
$ws.Range("O344").Value = 18000
$ws.Range("P344").Value = 18000
$ws.Range("Q344").Value = "$/caja 12 unidades"
$ws.Range("R344").Value = "Ecuador"
$ws.Range("S344").Value = 1500
$ws.Range("T344").Value = 12

# --- New row 345: Segunda, $/caja 14 unidades ---
$ws.Range("A345").Value = 3
$ws.Range("B345").Value = "Femacal de La Calera"
$ws.Range("C345").Value = "Coquimbo"
$ws.Range("D345").Value = 44551
$ws.Range("D345").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("E345").Value = 5
$ws.Range("F345").Value = "Fruta"
$ws.Range("G345").Value = 100108
$ws.Range("H345").Value = "Tropicales y subtropicales"
$ws.Range("I345").Value = 100108005
$ws.Range("J345").Value = "Piña"
$ws.Range("K345").Value = "Caramelo"
$ws.Range("L345").Value = "Segunda"
$ws.Range("M345").Value = 162
$ws.Range("N345").Value = 18000
$ws.Range("O345").Value = 18000
$ws.Range("P345").Value = 18000
$ws.Range("Q345").Value = "$/caja 14 unidades"
$ws.Range("R345").Value = "Ecuador"
$ws.Range("S345").Value = 1286
$ws.Range("T345").Value = 14
